# NMDC-EDGE Metagenomics ReadsQC bulk-submission template update
# - Update templates and handle empty rows and user input error

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the column header text (row 1)
#    C: "Paired-end Illumina/PacBio FASTQ" -> "Interleaved or Single-end Illumina/PacBio FASTQ"
#    D: "Illumina Pair-1 FASTQ"            -> "Illumina Paired-end R1 FASTQ"
#    E: "Illumina Pair-2 FASTQ"            -> "Illumina Paired-end  R2  FASTQ"
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = "Interleaved or Single-end Illumina/PacBio FASTQ"
$ws.Range("D1").Value = "Illumina Paired-end R1 FASTQ"
$ws.Range("E1").Value = "Illumina Paired-end  R2  FASTQ"

# ---------------------------------------------------------------------------
# 2. Remove the stray empty row 2 (it only carried a leftover wrap-text style
#    on C2, no real data) so the sheet goes back to a single header row plus
#    the data-entry rows below it.
# ---------------------------------------------------------------------------
$ws.Rows(2).Delete()

# ---------------------------------------------------------------------------
# 3. Resize columns to the new template layout
# ---------------------------------------------------------------------------
$ws.Columns(1).ColumnWidth = 28.5
$ws.Columns(2).ColumnWidth = 31.333333333333332
$ws.Columns(3).ColumnWidth = 40.166666666666664
$ws.Columns(4).ColumnWidth = 30.833333333333332
$ws.Columns(5).ColumnWidth = 30.833333333333332
$ws.Columns(6).ColumnWidth = 20.666666666666664

# ---------------------------------------------------------------------------
# 4. Rebuild the data validations (new prompts/titles, reordered so the
#    Illumina R1 FASTQ rule comes first) and handle the shrunk sqref ranges
#    (…2:…99 instead of …2:…100, already reflected by the row delete above).
# ---------------------------------------------------------------------------
$ws.Range("A2:A99").Validation.Delete()
$ws.Range("B2:B99").Validation.Delete()
$ws.Range("C2:C99").Validation.Delete()
$ws.Range("D2:D99").Validation.Delete()
$ws.Range("E2:E99").Validation.Delete()
$ws.Range("F2:F99").Validation.Delete()

# -- Illumina R1 FASTQ (D) -- now first in the list
$r = $ws.Range("D2:D99")
$r.Validation.Add(0, 1, 1)
$r.Validation.IgnoreBlank = $true
$r.Validation.InCellDropdown = $true
$r.Validation.ShowInput = $true
$r.Validation.ShowError = $true
$r.Validation.InputTitle = "Illumina R1 FASTQ"
$r.Validation.InputMessage = "Accept uploaded files, Retrieved SRA files and http(s) url inputs. Separate multiple files with commas.`n`nExamples:`nupload/test_R1.fq`nsra/SRR30724627_1.fastq.gz`nhttps://nmdc-edge.org/publicdata/test_data/Ecoli_10x.1.fastq"

# -- Project/Run Name (A)
$r = $ws.Range("A2:A99")
$r.Validation.Add(6, 1, 1, 3, 30)
$r.Validation.IgnoreBlank = $true
$r.Validation.InCellDropdown = $true
$r.Validation.ShowInput = $true
$r.Validation.ShowError = $true
$r.Validation.ErrorTitle = "Project/Run Name"
$r.Validation.ErrorMessage = "Invalid Input"
$r.Validation.InputTitle = "Project/Run Name"
$r.Validation.InputMessage = "Required. At least 3 but less than 30 characters. Only alphabets, numbers, dashs, dot and underscore are allowed."

# -- Description (B)
$r = $ws.Range("B2:B99")
$r.Validation.Add(0, 1, 1)
$r.Validation.IgnoreBlank = $true
$r.Validation.InCellDropdown = $true
$r.Validation.ShowInput = $true
$r.Validation.ShowError = $true
$r.Validation.InputTitle = "Description"
$r.Validation.InputMessage = "Optional"

# -- Illumina R2 FASTQ (E)
$r = $ws.Range("E2:E99")
$r.Validation.Add(0, 1, 1)
$r.Validation.IgnoreBlank = $true
$r.Validation.InCellDropdown = $true
$r.Validation.ShowInput = $true
$r.Validation.ShowError = $true
$r.Validation.InputTitle = "Illumina R2 FASTQ"
$r.Validation.InputMessage = "Accept uploaded files, Retrieved SRA files and http(s) url inputs. Separate multiple files with commas.`n`nExamples:`nupload/test_R2.fq`nsra/SRR30724627_2.fastq.gz`nhttps://nmdc-edge.org/publicdata/test_data/Ecoli_10x.2.fastq"

# -- Single Illumina/PacBio FASTQ (C)
$r = $ws.Range("C2:C99")
$r.Validation.Add(0, 1, 1)
$r.Validation.IgnoreBlank = $true
$r.Validation.InCellDropdown = $true
$r.Validation.ShowInput = $true
$r.Validation.ShowError = $true
$r.Validation.InputTitle = "Single Illumina/PacBio FASTQ"
$r.Validation.InputMessage = "Accept uploaded files, Retrieved SRA files and http(s) url inputs. Separate multiple files with commas.`n`nExamples:`nupload/test_R1.fq`nsra/SRR1602702.fastq.gz`nhttps://nmdc-edge.org/publicdata/test_data/Ecoli_interleaved_pe_small.fastq"
$r.Validation.ErrorTitle = ""
$r.Validation.ErrorMessage = ""

# -- Sequencing Platform (F)
$r = $ws.Range("F2:F99")
$r.Validation.Add(3, 1, 1, """Illumina, PacBio""")
$r.Validation.IgnoreBlank = $true
$r.Validation.InCellDropdown = $true
$r.Validation.ShowInput = $true
$r.Validation.ShowError = $true
$r.Validation.InputTitle = "Sequencing Platform"
$r.Validation.InputMessage = "Default: Illumina"
$r.Validation.ErrorTitle = ""
$r.Validation.ErrorMessage = ""

Write-Host "Template update applied."
